# Add a new "etat commande TMS" (TMS order status) column in Q, with a
# "valide" status for the existing order row, and select it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q1").Value = "etat commande TMS"
$ws.Range("Q2").Value = "valide"

$ws.Columns.Item(17).ColumnWidth = 26.3

$ws.Range("Q1:Q2").Select()
